# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates DAMSLTag (column I) and DialogAct (column J) values for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 10;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 14;  I = "%";  J = "Uninterpretable" },
    @{ Row = 23;  I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 27;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 44;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 45;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 48;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 62;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 69;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 78;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 82;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 83;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 86;  I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 89;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 91;  I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 101; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 106; I = "ba"; J = "Appreciation" },
    @{ Row = 108; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 109; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 130; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 131; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 136; I = "ba"; J = "Appreciation" },
    @{ Row = 137; I = "sv"; J = "Statement-opinion" },
    @{ Row = 152; I = "sv"; J = "Statement-opinion" },
    @{ Row = 158; I = "sv"; J = "Statement-opinion" },
    @{ Row = 168; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 178; I = "%";  J = "Uninterpretable" },
    @{ Row = 184; I = "sv"; J = "Statement-opinion" }
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 9).Value = $change.I
    $ws.Cells.Item($change.Row, 10).Value = $change.J
}
